$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 10 (shifts existing rows 10-19 down to 11-20)
$ws.Rows.Item(10).Insert()

# Populate the new row 10 with the same pattern as its neighboring rows,
# using the new date/volume values from the source update.
$ws.Cells.Item(10, 1).Value = 12
$ws.Cells.Item(10, 2).Value = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(10, 3).Value = "Metropolitana"
$ws.Cells.Item(10, 4).Value = 44497
$ws.Cells.Item(10, 4).NumberFormat = $ws.Cells.Item(11, 4).NumberFormat
$ws.Cells.Item(10, 5).Value = 13
$ws.Cells.Item(10, 6).Value = 100112028
$ws.Cells.Item(10, 7).Value = "Sandia"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 250
$ws.Cells.Item(10, 11).Value = 800
$ws.Cells.Item(10, 12).Value = 800
$ws.Cells.Item(10, 13).Value = 800
$ws.Cells.Item(10, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(10, 15).Value = "Perú"
$ws.Cells.Item(10, 16).Value = 800
$ws.Cells.Item(10, 17).Value = 1
$ws.Cells.Item(10, 18).Value = "Hortaliza"
